# Generate Report for Handback
# Updates row 6 ("9b1f07c5-6368-...md") on the zh-cn and de-de sheets with the
# latest handoff/handback info now that a target file + handback exist, plus
# widens the Error Detail column and records the "stale handback" message.

$wb = $excel.ActiveWorkbook

$sheetNames = @("zh-cn", "de-de")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Widen column P (Error Detail) to fit the new long message.
    $ws.Columns.Item(16).ColumnWidth = 39.17

    # I6: Latest Target File -> filename, with the same hyperlink styling used
    # by the other hyperlinked cells (A2:A6, I2:I5).
    $i6 = $ws.Range("I6")
    $i6.Value = "9b1f07c5-6368-4e48-a056-2d1b7b5f571a.md"
    $i6.Font.Underline = 2
    $i6.Font.Color = 15570276

    $target = "https://github.com/OpenLocalizationTestOrg/ol-test0-$($sheetName.Replace('-',''))/blob/93fb669b31b607e36023b90b379535d0194fab6e/e2e/9b1f07c5-6368-4e48-a056-2d1b7b5f571a.md"
    $ws.Hyperlinks.Add($i6, $target, "", "", "9b1f07c5-6368-4e48-a056-2d1b7b5f571a.md") | Out-Null

    # J6: Latest Handback File.
    $xlf = "9b1f07c5-6368-4e48-a056-2d1b7b5f571a.69e272cbbb7ed2a1207d4079bf4463274998d574.$sheetName.xlf"
    $ws.Range("J6").Value = $xlf

    # K6: Latest Handback DateTime.
    if ($sheetName -eq "zh-cn") {
        $ws.Range("K6").Value = "2016-11-15 17:02:55"
    } else {
        $ws.Range("K6").Value = "2016-11-15 17:03:13"
    }

    # P6: Error Detail - the handback isn't against the latest handoff yet.
    $ws.Range("P6").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/93fb669b31b607e36023b90b379535d0194fab6e/e2e/9b1f07c5-6368-4e48-a056-2d1b7b5f571a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/98c70c169f1fbf48888195390bde1f2753637b42/e2e/9b1f07c5-6368-4e48-a056-2d1b7b5f571a.md."
}

Write-Output "Handback report generated"
